$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new (blank) column before column N ("Late"), shifting the
# "Late" and "Outstanding" columns one place to the right to make room
# for a new column (used for variable instalments).
$ws.Columns("N:N").Insert()

# The newly inserted column gets an explicit width (not an auto/bestFit
# width like its neighbours).
$ws.Columns("N:N").ColumnWidth = 9.1

# Make "Repayment Schedule" the active sheet/tab and select cell R6 on it.
$ws.Activate()
$ws.Range("R6").Select()
